$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 949.0909
$ws.Range("I17").Value = 500
$ws.Range("J17").Value = 994
$ws.Range("K17").Value = 1500
$ws.Range("L17").Value = 2982
$ws.Range("M17").Value = -1332
$ws.Range("N17").Value = -3318

$ws.Range("H43").Value = 68578.47
$ws.Range("I43").Value = 91809.09
$ws.Range("J43").Value = 4694.25
$ws.Range("K43").Value = 91809.09
$ws.Range("L43").Value = 4694.25
$ws.Range("M43").Value = -91740.09
$ws.Range("N43").Value = -4832.25

$ws.Range("H133").Value = 38389.75
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 38389.75
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 38389.75
$ws.Range("N133").Value = -48509.75

$ws.Range("H137").Value = 17057998
$ws.Range("I137").Value = 6580152
$ws.Range("J137").Value = 25021162
$ws.Range("K137").Value = 19740456
$ws.Range("L137").Value = 75063486
$ws.Range("M137").Value = -19737906
$ws.Range("N137").Value = -75068586

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16483.83
$ws.Range("I32").Value = 13173.289
$ws.Range("J32").Value = 32647.059
$ws.Range("K32").Value = 13173.289
$ws.Range("L32").Value = 32647.059
$ws.Range("M32").Value = -12886.289
$ws.Range("N32").Value = -33221.059

$ws.Range("H74").Value = 13010133
$ws.Range("I74").Value = 919.2
$ws.Range("J74").Value = 33337028
$ws.Range("K74").Value = 919.2
$ws.Range("L74").Value = 33337028
$ws.Range("M74").Value = -45.20000000000005
$ws.Range("N74").Value = -33338776

$ws.Range("H77").Value = 13010133
$ws.Range("I77").Value = 919.2
$ws.Range("J77").Value = 33337028
$ws.Range("K77").Value = 4596
$ws.Range("L77").Value = 166685140
$ws.Range("M77").Value = -228
$ws.Range("N77").Value = -166693876

$ws.Range("H104").Value = 42980
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 42980
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 42980
$ws.Range("N104").Value = -49968

$ws.Range("H132").Value = 24932760
$ws.Range("I132").Value = 32963786
$ws.Range("J132").Value = 8335305
$ws.Range("K132").Value = 98891358
$ws.Range("L132").Value = 25005915
$ws.Range("M132").Value = -98888828
$ws.Range("N132").Value = -25010975

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5214266
$ws.Range("I31").Value = 9261029
$ws.Range("J31").Value = 11285.714
$ws.Range("K31").Value = 9261029
$ws.Range("L31").Value = 11285.714
$ws.Range("M31").Value = -9260734
$ws.Range("N31").Value = -11875.714

$ws.Range("H34").Value = 5214266
$ws.Range("I34").Value = 9261029
$ws.Range("J34").Value = 11285.714
$ws.Range("K34").Value = 9261029
$ws.Range("L34").Value = 11285.714
$ws.Range("M34").Value = -9260827
$ws.Range("N34").Value = -11689.714

$ws.Range("H122").Value = 3968.6177
$ws.Range("I122").Value = 4845.773
$ws.Range("J122").Value = 2360.5
$ws.Range("K122").Value = 14537.319
$ws.Range("L122").Value = 7081.5
$ws.Range("M122").Value = -12087.319
$ws.Range("N122").Value = -11981.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 248.21739
$ws.Range("I11").Value = 72.71429000000001
$ws.Range("J11").Value = 325
$ws.Range("K11").Value = 218.14287
$ws.Range("L11").Value = 975
$ws.Range("M11").Value = -78.14287000000002
$ws.Range("N11").Value = -1255

$ws.Range("H75").Value = 4501.625
$ws.Range("I75").Value = 1256.5
$ws.Range("J75").Value = 5583.3335
$ws.Range("K75").Value = 3769.5
$ws.Range("L75").Value = 16750.0005
$ws.Range("M75").Value = -2771.5
$ws.Range("N75").Value = -18746.0005

$ws.Range("H78").Value = 4501.625
$ws.Range("I78").Value = 1256.5
$ws.Range("J78").Value = 5583.3335
$ws.Range("K78").Value = 11308.5
$ws.Range("L78").Value = 50250.0015
$ws.Range("M78").Value = -6316.5
$ws.Range("N78").Value = -60234.0015

$ws.Range("H109").Value = 4781.067
$ws.Range("I109").Value = 1654.3334
$ws.Range("J109").Value = 6865.5557
$ws.Range("K109").Value = 4963.0002
$ws.Range("L109").Value = 20596.6671
$ws.Range("M109").Value = -3923.0002
$ws.Range("N109").Value = -22676.6671

$ws.Range("H121").Value = 3179426.2
$ws.Range("I121").Value = 554
$ws.Range("J121").Value = 4172824
$ws.Range("K121").Value = 1662
$ws.Range("L121").Value = 12518472
$ws.Range("M121").Value = -352
$ws.Range("N121").Value = -12521092

$ws.Range("H131").Value = 16245.877
$ws.Range("I131").Value = 71690.71000000001
$ws.Range("J131").Value = 1025.7255
$ws.Range("K131").Value = 215072.13
$ws.Range("L131").Value = 3077.1765
$ws.Range("M131").Value = -210032.13
$ws.Range("N131").Value = -13157.1765

$ws.Range("H132").Value = 1537.0588
$ws.Range("I132").Value = 745.3125
$ws.Range("J132").Value = 2240.8333
$ws.Range("K132").Value = 6707.8125
$ws.Range("L132").Value = 20167.4997
$ws.Range("M132").Value = -4177.8125
$ws.Range("N132").Value = -25227.4997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4373312
$ws.Range("I70").Value = 1789835
$ws.Range("J70").Value = 11908453
$ws.Range("K70").Value = 1789835
$ws.Range("L70").Value = 11908453
$ws.Range("M70").Value = -1789565
$ws.Range("N70").Value = -11908993

$ws.Range("H73").Value = 4373312
$ws.Range("I73").Value = 1789835
$ws.Range("J73").Value = 11908453
$ws.Range("K73").Value = 1789835
$ws.Range("L73").Value = 11908453
$ws.Range("M73").Value = -1788899
$ws.Range("N73").Value = -11910325

$ws.Range("H80").Value = 9644.375
$ws.Range("I80").Value = 4040.2632
$ws.Range("J80").Value = 30940
$ws.Range("K80").Value = 4040.2632
$ws.Range("L80").Value = 30940
$ws.Range("M80").Value = -3042.2632
$ws.Range("N80").Value = -32936

$ws.Range("H83").Value = 9644.375
$ws.Range("I83").Value = 4040.2632
$ws.Range("J83").Value = 30940
$ws.Range("K83").Value = 20201.316
$ws.Range("L83").Value = 154700
$ws.Range("M83").Value = -15209.316
$ws.Range("N83").Value = -164684

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1653.8334
$ws.Range("I7").Value = 1473.6471
$ws.Range("J7").Value = 2091.4285
$ws.Range("K7").Value = 1473.6471
$ws.Range("L7").Value = 2091.4285
$ws.Range("M7").Value = -1361.6471
$ws.Range("N7").Value = -2315.4285

$ws.Range("H40").Value = 3049.8333
$ws.Range("I40").Value = 3931.3333
$ws.Range("J40").Value = 2168.3333
$ws.Range("K40").Value = 3931.3333
$ws.Range("L40").Value = 2168.3333
$ws.Range("M40").Value = -3795.3333
$ws.Range("N40").Value = -2440.3333

$ws.Range("H126").Value = 1653.8334
$ws.Range("I126").Value = 1473.6471
$ws.Range("J126").Value = 2091.4285
$ws.Range("K126").Value = 4420.9413
$ws.Range("L126").Value = 6274.2855
$ws.Range("M126").Value = -1950.9413
$ws.Range("N126").Value = -11214.2855

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 42377
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 42377
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 42377
$ws.Range("N109").Value = -45151

$ws.Range("H132").Value = 1110591.9
$ws.Range("I132").Value = 3048.5676
$ws.Range("J132").Value = 7940442
$ws.Range("K132").Value = 9145.702799999999
$ws.Range("L132").Value = 23821326
$ws.Range("M132").Value = -6615.702799999999
$ws.Range("N132").Value = -23826386
